# Kumbham_LabExam03Grading.xlsx - "Driver class from 34-48"
#
# This script reproduces the grading edits made to the "CustomerMappingDriver
# Class" section (rows 28-31) of Sheet1:
#   - The comment for row 29 (addProduct/driver related row) is replaced with
#     a note about the missing scanner class / customer / product wiring.
#   - The comment for row 30 is replaced with a note about incorrect results
#     in the driver class, and its "Points for grading" (E30) earned score
#     drops from 8 to 7.
#   - The section subtotal (E31) and the grand total (E38) are formulas that
#     recalculate automatically once E29 changes.
#   - The active selection moves from F3 to F30, matching where the author
#     was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: CustomerMappingDriver Class, "For writing partial code for scanner class" ---
$ws.Range("F29").Value = "(-9) for not completing the scanner class, not declaring and initalizing customer and product object and adding them to inventory. Not writing else statement "

# Points earned for row 29 dropped from 8 to 7
$ws.Range("E29").Value = 7

# --- Row 30: CustomerMappingDriver Class, "For incorrect results for all methods" ---
$ws.Range("F30").Value = "(-4) For incorrect results for all methods in driver class"

# E31 (=SUM(E29:E30)) and E38 (grand total) recalculate automatically.

# Move the active selection to F30 (where the grading comment was last edited)
$ws.Range("F30").Select()
